# TC16_Canine_Filter_Breed-Chesapeake.xlsx - "corrected ICDC Breed 1-14 scripts"
#
# The "startup" sheet's B2/B4 cells hold Cypher query text used to drive the
# report. This edit:
#   1. Restores the Chinese Shar-Pei query (previously in B4) to B2.
#   2. Puts a corrected Chesapeake Bay Retriever "file" query into B4 with the
#      stray `File Type` and `Breed` output columns removed.
#   3. Shrinks row 4's height to match the now-shorter wrapped text.
#   4. Leaves the selection on B4, matching where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$sharPeiQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Chinese Shar-Pei']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
        coalesce(co.cohort_description, '') AS `Cohort`

'@

$chesapeakeFileQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Chesapeake Bay Retriever']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B2").Value = $sharPeiQuery
$ws.Range("B4").Value = $chesapeakeFileQuery

$ws.Rows.Item(4).RowHeight = 217.5

$ws.Range("B4").Select()
